# UC-47 "Editar natureza financeira" - Ajuste na descricao da pre-condicao
$d = $word.ActiveDocument

# 1) Pre-condicao bullet: era sobre "contratos", agora sobre "naturezas
#    financeiras" (e concordancia de genero: cadastrados -> cadastradas).
$r1 = $d.Content
$r1.Find.Execute(
    "Deve haver registros de contratos previamente cadastrados [Caso de Uso 45]",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Deve haver registros de naturezas financeiras previamente cadastradas [Caso de Uso 45]",
    2) | Out-Null

# 2) Ajuste de concordancia de genero na regra "codigo fixo".
$r2 = $d.Content
$r2.Find.Execute(
    " Uma vez cadastrado uma ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Uma vez cadastrada uma ",
    2) | Out-Null

# 3) O bookmark `_GoBack` (marcador de "ultima posicao editada" do Word)
#    precisa ser removido do local antigo (depois da imagem, ao final do
#    documento) e recriado logo depois do trecho que acabou de ser editado
#    em "...determinado codigo, o| mesmo nao podera...".
#
#    OBS: nao usar Range.Delete()/InsertXML() sobre um Range colapsado
#    (Start == End) - isso corrompe o documento neste runtime. Por isso a
#    remocao do bookmark antigo e feita reescrevendo, via InsertXML, um
#    Range de largura >= 1 que contem o bookmark (um caractere antes at
#    um caractere depois), preservando a imagem que esta no mesmo paragrafo.
$old = $d.Bookmarks.Item("_GoBack")
$oldPos = $old.Start
$wrapRange = $d.Range($oldPos - 1, $oldPos + 1)
$wrapRange.InsertXML(
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w:rsidR="00CA6C1B" w:rsidRDefault="0050201F" w:rsidP="00DD5DC9"><w:pPr><w:pStyle w:val="Cabealho"/><w:tabs><w:tab w:val="clear" w:pos="4252"/><w:tab w:val="clear" w:pos="8504"/><w:tab w:val="left" w:pos="395"/></w:tabs><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="808080" w:themeColor="background1" w:themeShade="80"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:noProof/><w:color w:val="808080" w:themeColor="background1" w:themeShade="80"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:pict><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1025" type="#_x0000_t75" style="width:453.75pt;height:515.25pt"><v:imagedata r:id="rId8" o:title="editar_natureza"/></v:shape></w:pict></w:r></w:p>'
) | Out-Null

# 4) Recria o bookmark `_GoBack` logo apos "...determinado codigo, o" e
#    antes de " mesmo nao podera ser alterado...".
$r3 = $d.Content
$r3.Find.Execute(
    "com determinado código, o",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$newPos = $r3.End
$newPoint = $d.Range($newPos, $newPos)
$d.Bookmarks.Add("_GoBack", $newPoint) | Out-Null
